$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (minor revisions to the series) ---
$ws.Range("C111").Value = 87.37
$ws.Range("C112").Value = 86.88
$ws.Range("C113").Value = 86.49
$ws.Range("C114").Value = 87.01
$ws.Range("C115").Value = 88.94
$ws.Range("C118").Value = 90.27
$ws.Range("C119").Value = 90.83
$ws.Range("C121").Value = 88.24
$ws.Range("C123").Value = 89.35
$ws.Range("C124").Value = 82.77
$ws.Range("C127").Value = 91.44
$ws.Range("C130").Value = 90.86
$ws.Range("C133").Value = 91.96
$ws.Range("C134").Value = 95.08
$ws.Range("C135").Value = 90.43
$ws.Range("C136").Value = 93.5
$ws.Range("C137").Value = 93.22
$ws.Range("C138").Value = 92.97
$ws.Range("C144").Value = 94.33
$ws.Range("C148").Value = 96.1
$ws.Range("C149").Value = 97.38
$ws.Range("C150").Value = 96.82
$ws.Range("C153").Value = 97.86
$ws.Range("C154").Value = 97.27
$ws.Range("C156").Value = 97.96
$ws.Range("C157").Value = 97.14
$ws.Range("C158").Value = 101.26
$ws.Range("C161").Value = 97.12
$ws.Range("C162").Value = 97.26
$ws.Range("C163").Value = 99.06
$ws.Range("C164").Value = 101.5
$ws.Range("C165").Value = 101.03
$ws.Range("C166").Value = 99.89
$ws.Range("C167").Value = 100.63
$ws.Range("C170").Value = 99.92
$ws.Range("C173").Value = 101.17
$ws.Range("C175").Value = 101.56
$ws.Range("C176").Value = 99.65
$ws.Range("C179").Value = 99.83
$ws.Range("C181").Value = 100.25
$ws.Range("C182").Value = 103.55
$ws.Range("C184").Value = 99.13
$ws.Range("C186").Value = 101.48
$ws.Range("C187").Value = 103.9
$ws.Range("C190").Value = 100.82
$ws.Range("C193").Value = 99.45
$ws.Range("C194").Value = 99.81
$ws.Range("C196").Value = 102.53
$ws.Range("C198").Value = 99.4
$ws.Range("C199").Value = 99.3
$ws.Range("C200").Value = 98.22
$ws.Range("C202").Value = 98.7
$ws.Range("C203").Value = 97.08
$ws.Range("C205").Value = 98.55
$ws.Range("C206").Value = 98.63
$ws.Range("C208").Value = 91.23
$ws.Range("C209").Value = 97.3
$ws.Range("C212").Value = 100.07
$ws.Range("C213").Value = 101.86
$ws.Range("C214").Value = 100.88
$ws.Range("C215").Value = 101.01
$ws.Range("C216").Value = 101.65
$ws.Range("C218").Value = 101.21
$ws.Range("C220").Value = 102.31
$ws.Range("C221").Value = 101.38
$ws.Range("C222").Value = 101.93
$ws.Range("C223").Value = 103.1
$ws.Range("C224").Value = 101.95
$ws.Range("C225").Value = 101.36
$ws.Range("C226").Value = 100.57
$ws.Range("C227").Value = 100.69
$ws.Range("C228").Value = 104.6
$ws.Range("C230").Value = 100.54
$ws.Range("C233").Value = 102.62
$ws.Range("C235").Value = 101.07
$ws.Range("C236").Value = 101.58
$ws.Range("C237").Value = 103.13
$ws.Range("C238").Value = 101.34
$ws.Range("C239").Value = 98.34
$ws.Range("C240").Value = 101.44
$ws.Range("C241").Value = 104.85
$ws.Range("C242").Value = 102.72
$ws.Range("C243").Value = 101.67
$ws.Range("C245").Value = 98.06
$ws.Range("C246").Value = 98.21
$ws.Range("C247").Value = 97.28
$ws.Range("C248").Value = 97.96
$ws.Range("C249").Value = 98.46
$ws.Range("B250").Value = 98.34
$ws.Range("C250").Value = 100.62
$ws.Range("C251").Value = 102.21
$ws.Range("C253").Value = 100.5
$ws.Range("C254").Value = 101.83
$ws.Range("C255").Value = 101.58
$ws.Range("C257").Value = 103.14
$ws.Range("C259").Value = 104.25
$ws.Range("C260").Value = 103.35
$ws.Range("B261").Value = 103.65
$ws.Range("C261").Value = 102.17

# --- Append new month row 262 (01-09-2021) ---
$ws.Range("A262").NumberFormat = "@"
$ws.Range("A262").Value = "01-09-2021"
$ws.Range("A262").ClearFormats()
$ws.Range("B262").Value = 97.61
$ws.Range("C262").Value = 100.08
